# Add plotRxnPathDiagram info into the Master_list.xlsx.
# A new entry is inserted as its own row right after the
# "mechanism_analyzer.ipynb" row (row 63) in the "Operations for many
# species" section, pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 64, shifting rows 64-81 down to 65-82.
$ws.Rows.Item(64).Insert() | Out-Null

$ws.Range("A64").Value = "plotRxnPathDiagram"
$ws.Range("B64").Value = "WIP/model_analyzer/plotRxnPathDiagram"
$ws.Range("C64").Value = "Peng"
$ws.Range("D64").Value = "Plot rxn path diagram for constant V simulation using Cantera. The nodes in the output diagram is labeled with species images."

# Match the author's final cursor position.
$ws.Range("D64").Select() | Out-Null
